# Applying odds updates for Jogos_da_Semana_FlashScore_2024-10-16.xlsx
# Commit: Atualizando o arquivo XLSX

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 1.73
$ws.Range("I3").Value = 5
$ws.Range("AC3").Value = 8
$ws.Range("AR3").Value = 51
$ws.Range("AZ3").Value = 126
$ws.Range("H5").Value = 3.9
$ws.Range("J5").Value = 2.2
$ws.Range("K5").Value = 2.1
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("W5").Value = 5.5
$ws.Range("X5").Value = 6.5
$ws.Range("AA5").Value = 15
$ws.Range("AB5").Value = 34
$ws.Range("AC5").Value = 8
$ws.Range("AD5").Value = 8
$ws.Range("AH5").Value = 11
$ws.Range("AJ5").Value = 19
$ws.Range("AM5").Value = 51
$ws.Range("AN5").Value = 3.4
$ws.Range("AO5").Value = 8
$ws.Range("AQ5").Value = 26
$ws.Range("AS5").Value = 201
$ws.Range("AX5").Value = 34
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 3.25
$ws.Range("K8").Value = 2.1
$ws.Range("L8").Value = 2.77
$ws.Range("N8").Value = 7.8
$ws.Range("O8").Value = 1.32
$ws.Range("P8").Value = 2.85
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.7
$ws.Range("U8").Value = 1.75
$ws.Range("V8").Value = 1.87
$ws.Range("W8").Value = 9.25
$ws.Range("X8").Value = 16
$ws.Range("AC8").Value = 9
$ws.Range("AE8").Value = 14.5
$ws.Range("AF8").Value = 70
$ws.Range("AG8").Value = 600
$ws.Range("AH8").Value = 7.2
$ws.Range("AJ8").Value = 9
$ws.Range("AL8").Value = 18.5
$ws.Range("AM8").Value = 30
$ws.Range("AN8").Value = 5
$ws.Range("AP8").Value = 23
$ws.Range("AR8").Value = 110
$ws.Range("AS8").Value = 300
$ws.Range("AT8").Value = 2.57
$ws.Range("AU8").Value = 6.9
$ws.Range("AV8").Value = 60
$ws.Range("AW8").Value = 4.05
$ws.Range("AX8").Value = 11.25
$ws.Range("AY8").Value = 19.5
$ws.Range("BA8").Value = 75
$ws.Range("BB8").Value = 250
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("G14").Value = 2.52
$ws.Range("I14").Value = 2.55
$ws.Range("J14").Value = 3.15
$ws.Range("K14").Value = 2.1
$ws.Range("L14").Value = 3.15
$ws.Range("O14").Value = 1.33
$ws.Range("P14").Value = 3.05
$ws.Range("Q14").Value = 2
$ws.Range("T14").Value = 2.72
$ws.Range("W14").Value = 8
$ws.Range("X14").Value = 12.5
$ws.Range("Y14").Value = 9.75
$ws.Range("Z14").Value = 28
$ws.Range("AA14").Value = 22
$ws.Range("AE14").Value = 14
$ws.Range("AI14").Value = 12.5
$ws.Range("AJ14").Value = 9.75
$ws.Range("AK14").Value = 28
$ws.Range("AL14").Value = 22
$ws.Range("AN14").Value = 4.5
$ws.Range("AO14").Value = 14
$ws.Range("AP14").Value = 22
$ws.Range("AQ14").Value = 60
$ws.Range("AR14").Value = 100
$ws.Range("AS14").Value = 300
$ws.Range("AT14").Value = 2.72
$ws.Range("AW14").Value = 4.5
$ws.Range("AX14").Value = 13.5
$ws.Range("AY14").Value = 21
$ws.Range("BA14").Value = 90
$ws.Range("BB14").Value = 250